$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 63
$ws.Range("A4").Value = 77
$ws.Range("A5").Value = 0
$ws.Range("A6").Value = 1

$ws.Range("A14").Value = 2
$ws.Range("A15").Value = 1
$ws.Range("A16").Value = 1
$ws.Range("A17").Value = 3

$ws.Range("A19").Value = 1
